$wb = $excel.ActiveWorkbook

# Rename sheets: Sheet1 -> GDN, Sheet3 -> ChiTiet
$wsGDN = $wb.Worksheets.Item("Sheet1")
$wsGDN.Name = "GDN"

$wsChiTiet = $wb.Worksheets.Item("Sheet3")
$wsChiTiet.Name = "ChiTiet"

# Update selection on GDN sheet (was H21, now D20) and clear tabSelected there
$wsGDN.Range("D20").Select()

# Activate ChiTiet sheet, making it the active/selected tab
$wsChiTiet.Activate()
